$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number need to be force-written as Text
# first (NumberFormat "@"), otherwise Excel auto-converts the literal string into a
# numeric value, which does not match the "price-as-text" cells used on this sheet.
# Style is reset to "Normal" afterwards so no stray explicit cell style is left behind.
$textForceCells = @(
    "D5",
    "D6",
    "D7",
    "D8",
    "D9",
    "D10",
    "D11",
    "D13",
    "D14",
    "D16",
    "D18",
    "D20",
    "D23",
    "D25",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D36",
    "D38",
    "D41",
    "D42",
    "D43",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '29.420.09'
$ws.Range("E2").Value = '  +0.16%  '
# Row 3
$ws.Range("D3").Value = '1.850.10'
$ws.Range("E3").Value = '  +0.21%  '
# Row 4
$ws.Range("E4").Value = '  +0.15%  '
# Row 5
$ws.Range("D5").Value = '240.77'
$ws.Range("E5").Value = '  +0.12%  '
# Row 6
$ws.Range("D6").Value = '0.6311'
$ws.Range("E6").Value = '  +0.07%  '
# Row 7
$ws.Range("D7").Value = '1.001'
$ws.Range("E7").Value = '  +0.08%  '
# Row 8
$ws.Range("D8").Value = '0.07714'
# Row 9
$ws.Range("D9").Value = '0.2943'
$ws.Range("E9").Value = '  -0.46%  '
# Row 10
$ws.Range("D10").Value = '24.55'
$ws.Range("E10").Value = '  +0.44%  '
# Row 11
$ws.Range("D11").Value = '0.07752'
$ws.Range("E11").Value = '  +0.66%  '
# Row 12
$ws.Range("D12").Value = '1.849.41'
$ws.Range("E12").Value = '  -0.73%  '
# Row 13
$ws.Range("D13").Value = '5.028'
$ws.Range("E13").Value = '  +0.72%  '
# Row 14
$ws.Range("D14").Value = '0.00001085'
$ws.Range("E14").Value = '  +8.45%  '
# Row 15
$ws.Range("E15").Value = '  -0.37%  '
# Row 16
$ws.Range("D16").Value = '83.72'
$ws.Range("E16").Value = '  +1.06%  '
# Row 17
$ws.Range("D17").Value = '2.099.85'
$ws.Range("E17").Value = '  -0.25%  '
# Row 18
$ws.Range("D18").Value = '6.155'
$ws.Range("E18").Value = '  +0.58%  '
# Row 19
$ws.Range("D19").Value = '29.448.45'
$ws.Range("E19").Value = '  +0.16%  '
# Row 20
$ws.Range("D20").Value = '229.58'
$ws.Range("E20").Value = '  +0.77%  '
# Row 21
$ws.Range("E21").Value = '  +0.34%  '
# Row 22
$ws.Range("E22").Value = '  +0.07%  '
# Row 23
$ws.Range("D23").Value = '7.459'
$ws.Range("E23").Value = '  -1.07%  '
# Row 24
$ws.Range("E24").Value = '  +0.10%  '
# Row 25
$ws.Range("D25").Value = '157.54'
$ws.Range("E25").Value = '  +0.25%  '
# Row 26
$ws.Range("D26").Value = '0.1391'
$ws.Range("E26").Value = '  -0.60%  '
# Row 27
$ws.Range("D27").Value = '8.359'
$ws.Range("E27").Value = '  +0.07%  '
# Row 28
$ws.Range("D28").Value = '17.69'
$ws.Range("E28").Value = '  +0.26%  '
# Row 29
$ws.Range("D29").Value = '1.471'
$ws.Range("E29").Value = '  +0.48%  '
# Row 30
$ws.Range("D30").Value = '1.310'
$ws.Range("E30").Value = '  +4.69%  '
# Row 31
$ws.Range("D31").Value = '0.05726'
$ws.Range("E31").Value = '  +0.95%  '
# Row 32
$ws.Range("D32").Value = '4.111'
$ws.Range("E32").Value = '  -0.27%  '
# Row 33
$ws.Range("D33").Value = '4.055'
$ws.Range("E33").Value = '  +0.83%  '
# Row 34
$ws.Range("D34").Value = '1.853'
$ws.Range("E34").Value = '  +0.60%  '
# Row 35
$ws.Range("E35").Value = '  +0.41%  '
# Row 36
$ws.Range("D36").Value = '0.7097'
$ws.Range("E36").Value = '  -0.89%  '
# Row 37
$ws.Range("E37").Value = '  -0.26%  '
# Row 38
$ws.Range("D38").Value = '2.781'
$ws.Range("E38").Value = '  +0.10%  '
# Row 39
$ws.Range("D39").Value = '1.230.24'
$ws.Range("E39").Value = '  -2.41%  '
# Row 40
$ws.Range("E40").Value = '  -0.57%  '
# Row 41
$ws.Range("D41").Value = '6.494'
$ws.Range("E41").Value = '  +4.29%  '
# Row 42
$ws.Range("D42").Value = '0.9141'
$ws.Range("E42").Value = '  +0.56%  '
# Row 43
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  +0.12%  '
# Row 44
$ws.Range("B44").Value = 'RocketPoolETH'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D44").Value = '2.008.85'
$ws.Range("E44").Value = '  -0.27%  '
# Row 45
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").Value = '101.58'
$ws.Range("E45").Value = '  +0.33%  '
# Row 46
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '66.34'
$ws.Range("E46").Value = '  +0.22%  '
# Row 47
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").Value = '0.00000000122'
$ws.Range("E47").Value = '  +5.00%  '
# Row 48
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").Value = '7.155'
$ws.Range("E48").Value = '  +1.37%  '
# Row 49
$ws.Range("B49").Value = 'TheSandbox'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D49").Value = '0.4016'
$ws.Range("E49").Value = '  -0.51%  '
# Row 50
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '9.037'
$ws.Range("E50").Value = '  -0.93%  '
# Row 51
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").Value = '1.687'
$ws.Range("E51").Value = '  +0.23%  '

# Reset style on the force-text cells to avoid leaving an explicit cell style behind
foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
